$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = '@'
$c.Value = '26.863.69'
$c.ClearFormats()
$c = $ws.Range("E2")
$c.NumberFormat = '@'
$c.Value = '  -1.75%  '
$c.ClearFormats()
$c = $ws.Range("D3")
$c.NumberFormat = '@'
$c.Value = '1.824.26'
$c.ClearFormats()
$c = $ws.Range("E3")
$c.NumberFormat = '@'
$c.Value = '  -1.78%  '
$c.ClearFormats()
$c = $ws.Range("D4")
$c.NumberFormat = '@'
$c.Value = '1.007'
$c.ClearFormats()
$c = $ws.Range("E4")
$c.NumberFormat = '@'
$c.Value = '  +0.62%  '
$c.ClearFormats()
$c = $ws.Range("D5")
$c.NumberFormat = '@'
$c.Value = '310.61'
$c.ClearFormats()
$c = $ws.Range("E5")
$c.NumberFormat = '@'
$c.Value = '  -1.16%  '
$c.ClearFormats()
$c = $ws.Range("E6")
$c.NumberFormat = '@'
$c.Value = '  +0.61%  '
$c.ClearFormats()
$c = $ws.Range("D7")
$c.NumberFormat = '@'
$c.Value = '0.4573'
$c.ClearFormats()
$c = $ws.Range("E7")
$c.NumberFormat = '@'
$c.Value = '  -0.92%  '
$c.ClearFormats()
$c = $ws.Range("D8")
$c.NumberFormat = '@'
$c.Value = '0.3673'
$c.ClearFormats()
$c = $ws.Range("E8")
$c.NumberFormat = '@'
$c.Value = '  -1.06%  '
$c.ClearFormats()
$c = $ws.Range("D9")
$c.NumberFormat = '@'
$c.Value = '0.07149'
$c.ClearFormats()
$c = $ws.Range("D10")
$c.NumberFormat = '@'
$c.Value = '0.8716'
$c.ClearFormats()
$c = $ws.Range("E10")
$c.NumberFormat = '@'
$c.Value = '  -1.02%  '
$c.ClearFormats()
$c = $ws.Range("D11")
$c.NumberFormat = '@'
$c.Value = '0.07781'
$c.ClearFormats()
$c = $ws.Range("E11")
$c.NumberFormat = '@'
$c.Value = '  -0.09%  '
$c.ClearFormats()
$c = $ws.Range("D12")
$c.NumberFormat = '@'
$c.Value = '19.49'
$c.ClearFormats()
$c = $ws.Range("E12")
$c.NumberFormat = '@'
$c.Value = '  -2.02%  '
$c.ClearFormats()
$c = $ws.Range("D13")
$c.NumberFormat = '@'
$c.Value = '1.791.29'
$c.ClearFormats()
$c = $ws.Range("E13")
$c.NumberFormat = '@'
$c.Value = '  -3.26%  '
$c.ClearFormats()
$c = $ws.Range("D14")
$c.NumberFormat = '@'
$c.Value = '5.309'
$c.ClearFormats()
$c = $ws.Range("E14")
$c.NumberFormat = '@'
$c.Value = '  -1.43%  '
$c.ClearFormats()
$c = $ws.Range("D15")
$c.NumberFormat = '@'
$c.Value = '6.372'
$c.ClearFormats()
$c = $ws.Range("E15")
$c.NumberFormat = '@'
$c.Value = '  -2.68%  '
$c.ClearFormats()
$c = $ws.Range("D16")
$c.NumberFormat = '@'
$c.Value = '86.76'
$c.ClearFormats()
$c = $ws.Range("E16")
$c.NumberFormat = '@'
$c.Value = '  -5.52%  '
$c.ClearFormats()
$c = $ws.Range("D17")
$c.NumberFormat = '@'
$c.Value = '1.008'
$c.ClearFormats()
$c = $ws.Range("E17")
$c.NumberFormat = '@'
$c.Value = '  +0.78%  '
$c.ClearFormats()
$c = $ws.Range("D18")
$c.NumberFormat = '@'
$c.Value = '0.000008685'
$c.ClearFormats()
$c = $ws.Range("E18")
$c.NumberFormat = '@'
$c.Value = '  -4.50%  '
$c.ClearFormats()
$c = $ws.Range("D19")
$c.NumberFormat = '@'
$c.Value = '1.006'
$c.ClearFormats()
$c = $ws.Range("E19")
$c.NumberFormat = '@'
$c.Value = '  +0.51%  '
$c.ClearFormats()
$c = $ws.Range("D20")
$c.NumberFormat = '@'
$c.Value = '26.880.25'
$c.ClearFormats()
$c = $ws.Range("E20")
$c.NumberFormat = '@'
$c.Value = '  -1.73%  '
$c.ClearFormats()
$c = $ws.Range("D21")
$c.NumberFormat = '@'
$c.Value = '14.41'
$c.ClearFormats()
$c = $ws.Range("E21")
$c.NumberFormat = '@'
$c.Value = '  -2.51%  '
$c.ClearFormats()
$c = $ws.Range("D22")
$c.NumberFormat = '@'
$c.Value = '4.984'
$c.ClearFormats()
$c = $ws.Range("E22")
$c.NumberFormat = '@'
$c.Value = '  -2.81%  '
$c.ClearFormats()
$c = $ws.Range("D23")
$c.NumberFormat = '@'
$c.Value = '2.046.34'
$c.ClearFormats()
$c = $ws.Range("E23")
$c.NumberFormat = '@'
$c.Value = '  -4.38%  '
$c.ClearFormats()
$c = $ws.Range("D24")
$c.NumberFormat = '@'
$c.Value = '10.43'
$c.ClearFormats()
$c = $ws.Range("E24")
$c.NumberFormat = '@'
$c.Value = '  -0.76%  '
$c.ClearFormats()
$c = $ws.Range("D25")
$c.NumberFormat = '@'
$c.Value = '2.009'
$c.ClearFormats()
$c = $ws.Range("E25")
$c.NumberFormat = '@'
$c.Value = '  +4.09%  '
$c.ClearFormats()
$c = $ws.Range("D26")
$c.NumberFormat = '@'
$c.Value = '150.95'
$c.ClearFormats()
$c = $ws.Range("E26")
$c.NumberFormat = '@'
$c.Value = '  -0.83%  '
$c.ClearFormats()
$c = $ws.Range("D27")
$c.NumberFormat = '@'
$c.Value = '18.11'
$c.ClearFormats()
$c = $ws.Range("E27")
$c.NumberFormat = '@'
$c.Value = '  -1.32%  '
$c.ClearFormats()
$c = $ws.Range("D28")
$c.NumberFormat = '@'
$c.Value = '1.945'
$c.ClearFormats()
$c = $ws.Range("E28")
$c.NumberFormat = '@'
$c.Value = '  -6.06%  '
$c.ClearFormats()
$c = $ws.Range("D29")
$c.NumberFormat = '@'
$c.Value = '113.36'
$c.ClearFormats()
$c = $ws.Range("E29")
$c.NumberFormat = '@'
$c.Value = '  -2.29%  '
$c.ClearFormats()
$c = $ws.Range("D30")
$c.NumberFormat = '@'
$c.Value = '4.899'
$c.ClearFormats()
$c = $ws.Range("E30")
$c.NumberFormat = '@'
$c.Value = '  -4.00%  '
$c.ClearFormats()
$c = $ws.Range("D31")
$c.NumberFormat = '@'
$c.Value = '0.08797'
$c.ClearFormats()
$c = $ws.Range("E31")
$c.NumberFormat = '@'
$c.Value = '  -0.70%  '
$c.ClearFormats()
$c = $ws.Range("D32")
$c.NumberFormat = '@'
$c.Value = '2.996'
$c.ClearFormats()
$c = $ws.Range("E32")
$c.NumberFormat = '@'
$c.Value = '  -1.42%  '
$c.ClearFormats()
$c = $ws.Range("D33")
$c.NumberFormat = '@'
$c.Value = '0.7464'
$c.ClearFormats()
$c = $ws.Range("E33")
$c.NumberFormat = '@'
$c.Value = '  -3.38%  '
$c.ClearFormats()
$c = $ws.Range("D34")
$c.NumberFormat = '@'
$c.Value = '4.463'
$c.ClearFormats()
$c = $ws.Range("E34")
$c.NumberFormat = '@'
$c.Value = '  -0.70%  '
$c.ClearFormats()
$c = $ws.Range("D35")
$c.NumberFormat = '@'
$c.Value = '1.128'
$c.ClearFormats()
$c = $ws.Range("E35")
$c.NumberFormat = '@'
$c.Value = '  -3.83%  '
$c.ClearFormats()
$c = $ws.Range("D36")
$c.NumberFormat = '@'
$c.Value = '2.526'
$c.ClearFormats()
$c = $ws.Range("E36")
$c.NumberFormat = '@'
$c.Value = '  -4.85%  '
$c.ClearFormats()
$c = $ws.Range("E37")
$c.NumberFormat = '@'
$c.Value = '  +0.53%  '
$c.ClearFormats()
$c = $ws.Range("D38")
$c.NumberFormat = '@'
$c.Value = '0.01933'
$c.ClearFormats()
$c = $ws.Range("E38")
$c.NumberFormat = '@'
$c.Value = '  -1.31%  '
$c.ClearFormats()
$c = $ws.Range("D39")
$c.NumberFormat = '@'
$c.Value = '2.916'
$c.ClearFormats()
$c = $ws.Range("E39")
$c.NumberFormat = '@'
$c.Value = '  -1.19%  '
$c.ClearFormats()
$c = $ws.Range("D40")
$c.NumberFormat = '@'
$c.Value = '0.05103'
$c.ClearFormats()
$c = $ws.Range("E40")
$c.NumberFormat = '@'
$c.Value = '  -2.37%  '
$c.ClearFormats()
$c = $ws.Range("D41")
$c.NumberFormat = '@'
$c.Value = '6.910'
$c.ClearFormats()
$c = $ws.Range("E41")
$c.NumberFormat = '@'
$c.Value = '  -1.49%  '
$c.ClearFormats()
$c = $ws.Range("D42")
$c.NumberFormat = '@'
$c.Value = '0.4954'
$c.ClearFormats()
$c = $ws.Range("E42")
$c.NumberFormat = '@'
$c.Value = '  -3.71%  '
$c.ClearFormats()
$c = $ws.Range("E43")
$c.NumberFormat = '@'
$c.Value = '  -2.76%  '
$c.ClearFormats()
$c = $ws.Range("D44")
$c.NumberFormat = '@'
$c.Value = '8.260'
$c.ClearFormats()
$c = $ws.Range("E44")
$c.NumberFormat = '@'
$c.Value = '  -1.70%  '
$c.ClearFormats()
$c = $ws.Range("D45")
$c.NumberFormat = '@'
$c.Value = '0.4663'
$c.ClearFormats()
$c = $ws.Range("E45")
$c.NumberFormat = '@'
$c.Value = '  -3.37%  '
$c.ClearFormats()
$c = $ws.Range("D46")
$c.NumberFormat = '@'
$c.Value = '1.007'
$c.ClearFormats()
$c = $ws.Range("E46")
$c.NumberFormat = '@'
$c.Value = '  +0.65%  '
$c.ClearFormats()
$c = $ws.Range("D47")
$c.NumberFormat = '@'
$c.Value = '10.10'
$c.ClearFormats()
$c = $ws.Range("E47")
$c.NumberFormat = '@'
$c.Value = '  -2.23%  '
$c.ClearFormats()
$c = $ws.Range("D48")
$c.NumberFormat = '@'
$c.Value = '101.23'
$c.ClearFormats()
$c = $ws.Range("E48")
$c.NumberFormat = '@'
$c.Value = '  -1.67%  '
$c.ClearFormats()
$c = $ws.Range("D49")
$c.NumberFormat = '@'
$c.Value = '1.605'
$c.ClearFormats()
$c = $ws.Range("E49")
$c.NumberFormat = '@'
$c.Value = '  -2.82%  '
$c.ClearFormats()
$c = $ws.Range("D50")
$c.NumberFormat = '@'
$c.Value = '0.06091'
$c.ClearFormats()
$c = $ws.Range("E50")
$c.NumberFormat = '@'
$c.Value = '  -2.03%  '
$c.ClearFormats()
$c = $ws.Range("D51")
$c.NumberFormat = '@'
$c.Value = '64.33'
$c.ClearFormats()
$c = $ws.Range("E51")
$c.NumberFormat = '@'
$c.Value = '  -2.01%  '
$c.ClearFormats()
